# Remove the frequency-related rows from the test automation configuration.
# Rows 19-26 (GnssPoseSimulink, PointsRawFloat32, ImageRaw, ClockFrequency,
# SimulinkState, CurrentVelocity, PoseOtherCar, CurrentPose) are deleted
# entirely; this shifts the following rows (percent_reflecting_sfc, R) up
# so they become rows 19 and 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Range("A19:B26").EntireRow.Select()
$ws.Range("A19:B26").EntireRow.Delete()
